# Swap the presentation's theme palette from "Integral" to the stock
# "Office Theme" palette (ppt/theme/theme1.xml, the theme used by the
# slide master / every slide), matching the colours that used to live in
# ppt/theme/theme2.xml (the notes-master-only theme).
#
# clrScheme/fontScheme/fmtScheme are otherwise byte-identical between the
# two themes in this deck, so only the 12 scheme colours need to change.

function ToBgrInt([int]$r, [int]$g, [int]$b) {
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# index -> (name, target "Office Theme" RGB hex)
$officeTheme = @(
    @{ Index = 1;  Name = "dk1";      R = 0x00; G = 0x00; B = 0x00 },
    @{ Index = 2;  Name = "lt1";      R = 0xFF; G = 0xFF; B = 0xFF },
    @{ Index = 3;  Name = "dk2";      R = 0x44; G = 0x54; B = 0x6A },
    @{ Index = 4;  Name = "lt2";      R = 0xE7; G = 0xE6; B = 0xE6 },
    @{ Index = 5;  Name = "accent1";  R = 0x5B; G = 0x9B; B = 0xD5 },
    @{ Index = 6;  Name = "accent2";  R = 0xED; G = 0x7D; B = 0x31 },
    @{ Index = 7;  Name = "accent3";  R = 0xA5; G = 0xA5; B = 0xA5 },
    @{ Index = 8;  Name = "accent4";  R = 0xFF; G = 0xC0; B = 0x00 },
    @{ Index = 9;  Name = "accent5";  R = 0x44; G = 0x72; B = 0xC4 },
    @{ Index = 10; Name = "accent6";  R = 0x70; G = 0xAD; B = 0x47 },
    @{ Index = 11; Name = "hlink";    R = 0x05; G = 0x63; B = 0xC1 },
    @{ Index = 12; Name = "folHlink"; R = 0x95; G = 0x4F; B = 0x72 }
)

foreach ($entry in $officeTheme) {
    $color = $tcs.Colors($entry.Index)
    $color.RGB = ToBgrInt $entry.R $entry.G $entry.B
}

Write-Output "Theme colours updated to Office Theme palette."
